$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data record needs to be inserted as row 283, pushing all
# existing rows from 283 downward by one (old row 283 becomes row 284, etc).
$ws.Rows.Item(283).Insert()

# The non-varying columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical to the
# row that is now directly below the inserted one (row 284), so copy that
# row's formatting/values down into the new row first...
$ws.Rows.Item(284).Copy()
$ws.Rows.Item(283).PasteSpecial()
$excel.CutCopyMode = $false

# ...then overwrite the columns that differ for this new record.
$ws.Range("D283").Value = 44736
$ws.Range("J283").Value = 250
$ws.Range("K283").Value = 3100
$ws.Range("L283").Value = 3300
$ws.Range("M283").Value = 3196
$ws.Range("P283").Value = 533
